$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 5) with the latest log entry.
# Force text format on the date/time-looking columns so Excel does not
# auto-convert them into date/time serial values.
$ws.Range("C5:D5").NumberFormat = "@"

$ws.Range("A5").Value = "197cc1786f88b471"
$ws.Range("B5").Value = "thang truong <truongthoithang@gmail.com>"
$ws.Range("C5").Value = "2025-07-02"
$ws.Range("D5").Value = "17:02:51"
$ws.Range("E5").Value = "0 days, 0 hours, 0 minutes"
$ws.Range("F5").Value = "Not started"
